$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror the formatting of the existing header cells (e.g. G1) onto the
# new "Save" header cell so it gets the same style index instead of a
# freshly minted one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for the added column (plain, unstyled numeric cell).
$ws.Range("H2").Value = 0
